$d = $word.ActiveDocument

$map = @{
  "391×8=" = "185×7="
  "661×2=" = "829×8="
  "930×8=" = "391×4="
  "292×2=" = "973×7="
  "568×8=" = "832×9="
  "688×3=" = "177×7="
  "862×4=" = "903×9="
  "568×4=" = "408×2="
  "384×7=" = "466×5="
  "870×3=" = "580×2="
  "450×9=" = "982×3="
  "869×4=" = "285×5="
  "639×7=" = "602×7="
  "544×5=" = "927×5="
  "281×7=" = "748×9="
  "752×7=" = "383×2="
  "330×2=" = "386×8="
  "693×4=" = "563×5="
  "880×4=" = "245×2="
  "948×8=" = "792×4="
  "741×5=" = "675×7="
  "287×4=" = "211×7="
  "123×2=" = "289×7="
  "829×4=" = "542×2="
  "230×4=" = "374×4="
}

foreach ($old in $map.Keys) {
  $new = $map[$old]
  $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                           $true, 1, $false, $new, 2)
}
